$wb = $excel.ActiveWorkbook

# --- Rename sheet 2 from "Ebay" to "ExpediaTestData" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ExpediaTestData"

# --- Populate the new test-data table on ExpediaTestData ---
# Values are written in this specific order so that the resulting shared-string
# table is built up in the same sequence as the target workbook.
$ws2.Range("A1").Value = "PickUp"
$ws2.Range("B1").Value = "DropOff"
$ws2.Range("A2").Value = "LGA"
$ws2.Range("G2").Value = "January"
$ws2.Range("C1").Value = "PickUpDay"
$ws2.Range("A3").Value = "Newark"
$ws2.Range("D3").Value = "February"
$ws2.Range("D1").Value = "PickMonth"
$ws2.Range("E1").Value = "PickYear"
$ws2.Range("F1").Value = "DropDay"
$ws2.Range("G1").Value = "DropMonth"
$ws2.Range("G3").Value = "March"
# Leading apostrophe -> stored as text with the quotePrefix style (numeric-looking values)
$ws2.Range("C3").Value = "'19"
$ws2.Range("E3").Value = "'2022"
$ws2.Range("F2").Value = "'11"
$ws2.Range("F3").Value = "'5"
$ws2.Range("E2").Value = "'2021"
$ws2.Range("D2").Value = "December"
$ws2.Range("C2").Value = "'4"

$ws2.Range("B2").Value = "LGA"
$ws2.Range("B3").Value = "Newark"

# --- Column widths for the new table ---
# (ColumnWidth values chosen so the saved OOXML "width" attribute comes out to
# exactly 14 / 15 / 12.5, matching Excel's internal character-width rounding.)
$ws2.Columns.Item(3).ColumnWidth = 13.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 14.166666666666668
$ws2.Columns.Item(7).ColumnWidth = 11.666666666666668

# --- Make ExpediaTestData the active/selected tab, with E9 selected ---
$ws2.Activate() | Out-Null
$ws2.Range("E9").Select() | Out-Null
